# Apply mission file parameter updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Parameter" column (C) for rows 4-9 to 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 1

# Update the active selection to E9
$ws.Range("E9").Select()
